$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values after re-pulling data / mean calculation
$updates = @{
    3  = -2
    6  = -3
    7  = 0
    12 = 0
    13 = 2
    15 = 1
    22 = -1
    24 = -1
    27 = 0
    28 = 1
    32 = 1
    33 = 0
    34 = -4
    35 = 2
    36 = 0
    40 = 2
    43 = 2
    47 = -2
    54 = 5
    55 = 1
    62 = 0
    68 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
